$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("deposits")
$ws1.Range("B6").Value = "AAPL"

$ws2 = $wb.Worksheets.Item("dividends")
$ws2.Range("B3").Value = "AAPL"

$ws3 = $wb.Worksheets.Item("sales")
$ws3.Range("B6").Value = "AAPL"
$ws3.Range("B7").Value = "AAPL"
